$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "1" markers at A2, B3, C5, D5 with the text value "0.5" without
# disturbing the cells' existing (General) number format/style: build the text
# in a scratch cell via a formula (so it is typed as text, not auto-coerced to
# a number), then paste just the value into the target and clean up.
foreach ($addr in @("A2", "B3", "C5", "D5")) {
    $ws.Range("F1").Formula = '="0.5"'
    $ws.Range("F1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null  # xlPasteValues
}
$ws.Range("F1").Clear() | Out-Null

# The old "1" at C4 becomes 0.
$ws.Range("C4").Value = 0

# Move the active selection to D6 as in the edited workbook.
$ws.Range("D6").Select() | Out-Null
